# august 13 catch update
# - Update row 28 (FN0774) Notes to clarify selective-gear-only requirement
# - Add new row 34 for FN0809 (Aug 13-19 sockeye opening)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skeena")

# New catch-update row for FN0809 (Aug 13-19)
$ws.Range("A34").Value = "FN0809"
$ws.Range("B34").Value = "Aboriginal"
$ws.Range("C34").Value = "Aug 13-19"
$ws.Range("D34").Value = "Sockeye"
$ws.Range("E34").Value = "Selective Gear"
$ws.Range("F34").Value = "Region 6-Gitksan"
$ws.Range("G34").Value = 7
$ws.Range("I34").Value = "Sockeye target, selective gear only"

# Match the date-column formatting used by the other "Aug ##-##" rows
$ws.Range("C28").Copy()
$ws.Range("C34").PasteSpecial(-4122)

# Clarify the existing FN0774 notes to match the new wording
$ws.Range("I28").Value = "Sockeye target, selective gear only"

# Update the on-screen selection to reflect where the edit was made
$ws.Range("I29").Select()
